$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the micro-precision of the "Fecha" timestamp for the existing
#    last batch of rows (730-743) - same displayed date/time, tiny float fix.
$fixedTimestamp = 44232.76784993055
for ($r = 730; $r -le 743; $r++) {
    $ws.Range("D" + $r).Value = $fixedTimestamp
}

# 2) Append a brand-new batch of 14 rows (744-757) - one row per monitored
#    service, mirroring the structure of every earlier batch in the sheet.
#    Columns: A = Nombre, B = URL (hyperlink), C = Disponibilidad, D = Fecha
$newTimestamp = 44232.78916576152

# name, hyperlink target address, sub-address/location, displayed text
$rows = @(
    @("Odoo",              "https://www.dataintelligence-group.com/",                    "", "https://www.dataintelligence-group.com/"),
    @("Blackbox",          "https://serviciodashboard.azurewebsites.net/",                "", "https://serviciodashboard.azurewebsites.net/"),
    @("PowerBI",           "https://powerbi.microsoft.com/es-es/",                        "", "https://powerbi.microsoft.com/es-es/"),
    @("Dropbox",           "https://www.dropbox.com/",                                    "", "https://www.dropbox.com/"),
    @("Odoo",              "https://dataintelligence.store/",                             "", "https://dataintelligence.store/"),
    @("GEE",                "https://app-data-i.users.earthengine.app/",                  "", "https://app-data-i.users.earthengine.app/"),
    @("UtilidadesOdoo",    "https://odooutil.azurewebsites.net/",                         "", "https://odooutil.azurewebsites.net/"),
    @("Filtros Dashboard", "https://filtradordashboard.azurewebsites.net/",               "", "https://filtradordashboard.azurewebsites.net/"),
    @("MapStore",          "https://ide.dataintelligence-group.com/mapstore/",            "/", "https://ide.dataintelligence-group.com/mapstore/#/"),
    @("GeoServer",         "https://ide.dataintelligence-group.com/geoserver/web/?0",     "", "https://ide.dataintelligence-group.com/geoserver/web/?0"),
    @("Tomcat",            "https://ide.dataintelligence-group.com/",                     "", "https://ide.dataintelligence-group.com/"),
    @("Shiny",             "https://rpubs.com/dataintelligence/",                         "", "https://rpubs.com/dataintelligence/"),
    @("Github",            "https://github.com/Sud-Austral/",                             "", "https://github.com/Sud-Austral/"),
    @("EZ Exporter",       "https://ezexporter.highviewapps.com/exports/export-profile/", "", "https://ezexporter.highviewapps.com/exports/export-profile/")
)

$row = 744
foreach ($item in $rows) {
    $name = $item[0]
    $url = $item[1]
    $location = $item[2]
    $displayText = $item[3]

    $ws.Range("A" + $row).Value = $name
    $ws.Range("C" + $row).Value = "Disponible"

    if ($location -ne "") {
        $ws.Hyperlinks.Add($ws.Range("B" + $row), $url, $location)
    } else {
        $ws.Hyperlinks.Add($ws.Range("B" + $row), $url)
    }
    # Keep the displayed cell text pinned (reuses existing shared strings,
    # e.g. the MapStore URL keeps its "#/" fragment in the visible text while
    # the actual relationship target is the bare URL + location).
    $ws.Range("B" + $row).Value = $displayText
    $ws.Range("B" + $row).Style = "Hyperlink"

    $ws.Range("D" + $row).Value = $newTimestamp
    $ws.Range("D" + $row).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $row = $row + 1
}
